# The deck ships two theme parts: ppt/theme/theme1.xml ("Office Theme" /
# stock blue-and-orange palette) and ppt/theme/theme2.xml ("Integral" /
# Red-Violet palette). theme2.xml is the theme actually wired to the
# slide master (and therefore to every slide), so it's the one whose
# colours are visible throughout the deck. The edit swaps the deck back
# to the plain "Office" colour palette by recolouring the active theme's
# colour scheme, one slot at a time, through the SlideMaster.

$p = $ppt.ActivePresentation
$cs = $p.SlideMaster.ColorScheme

# Office colour scheme (Background/Text dark-1/light-1/dark-2/light-2,
# Accent 1-6, Hyperlink, Followed Hyperlink) expressed as COM RGB longs
# (0xBBGGRR) - the same values PowerPoint itself writes for the stock
# "Office" theme palette.
$cs.Item(1).RGB  = 0         # Text/Background Dark 1   -> 000000
$cs.Item(2).RGB  = 16777215  # Text/Background Light 1  -> FFFFFF
$cs.Item(3).RGB  = 6968388   # Text/Background Dark 2   -> 44546A
$cs.Item(4).RGB  = 15132391  # Text/Background Light 2  -> E7E6E6
$cs.Item(5).RGB  = 13998939  # Accent 1                 -> 5B9BD5
$cs.Item(6).RGB  = 3243501   # Accent 2                 -> ED7D31
$cs.Item(7).RGB  = 10855845  # Accent 3                 -> A5A5A5
$cs.Item(8).RGB  = 49407     # Accent 4                 -> FFC000
$cs.Item(9).RGB  = 12874308  # Accent 5                 -> 4472C4
$cs.Item(10).RGB = 4697456   # Accent 6                 -> 70AD47
$cs.Item(11).RGB = 12673797  # Hyperlink                -> 0563C1
$cs.Item(12).RGB = 7491477   # Followed Hyperlink       -> 954F72
